$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Search_Request")

$ws.Range("A2").Value = 150286
$ws.Range("A3").Value = 150286
$ws.Range("A4").Value = 122075
$ws.Range("A5").Value = 108054
